$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed during the repull/push of data
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -4
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -4
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = -3
